$d = $word.ActiveDocument

# 1. Update the version/date line.
$d.Content.Find.Execute(
    "Wersja: 1.0  •  Data: 2026-02-06 17:27  •  Obszar: dom + warsztat (≈2000 m²)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Wersja: 1.0  •  Data: 2026-02-07 17:08  •  Obszar: dom + warsztat (≈2000 m²)",
    2
)

# 2. Insert a new row into the "Załącznik A" parameters table, right before
#    the row describing pumpRestoreDryMinutes, adding pumpOffRepeatMinSeconds.

# Find the parameters table: the one whose first cell of row 1 reads "Parametr".
$table = $null
for ($t = 1; $t -le $d.Tables.Count; $t++) {
    $candidate = $d.Tables.Item($t)
    if ($candidate.Cell(1, 1).Range.Text.TrimEnd([char]7, [char]13) -eq "Parametr") {
        $table = $candidate
        break
    }
}

$targetRow = $null
for ($i = 1; $i -le $table.Rows.Count; $i++) {
    if ($table.Cell($i, 1).Range.Text.TrimEnd([char]7, [char]13) -eq "pumpRestoreDryMinutes") {
        $targetRow = $i
        break
    }
}

$newRow = $table.Rows.Add($table.Rows.Item($targetRow))
$newRow.Cells.Item(1).Range.Text = "pumpOffRepeatMinSeconds"
$newRow.Cells.Item(2).Range.Text = "60"
$newRow.Cells.Item(3).Range.Text = "Minimalny odstęp ponawiania komendy OFF dla pomp podczas aktywnego zalania (jeśli pierwsza komenda zginęła w obciążeniu)."
